$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    3 = @(0.6606524410359556, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.611132179096228)
    4 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    5 = @(1.455362044514542, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 2.405841782574814)
    6 = @(0.0006408296065709695, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 4.339460241733758)
    7 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    8 = @(3.286832544864788, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 27.82738278199502)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
